$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Elemental Shaman', ['Token Creature — Elemental Shaman', '3/1'])"
$ws.Range("A3:A4").ClearContents()
